# 2428-RBI-EI-DB-DL-REC-NOCOM-RNI-CTPD-DL-MD-TR-1-ADD-VAR-INST-Newcreateloan1.xlsx
# "Loan RBI, Variable Instalments"
#
# The "Repayment schedule" sheet gets a new (empty) column inserted before
# the existing "Late" column, shifting Late / Outstanding-heading /
# Outstanding-value one column to the right (N -> O -> P -> Q), to make room
# for a variable-instalment "Late" style amount column used by the new
# RBI variable-instalment test flow. The new column keeps the width of the
# column immediately to its left (column M), matching what Excel does when
# a column is inserted via Insert Cells/Columns.
#
# In addition, the active sheet/selection bookmarks change: "Repayment
# schedule" becomes the active tab (selection resting on N18) instead of
# "Edit Repayment Schedule" (whose selection moves to C7, no longer the
# active tab).

$wb = $excel.ActiveWorkbook

$repayment = $wb.Worksheets.Item("Repayment schedule")
$editSchedule = $wb.Worksheets.Item("Edit Repayment Schedule")

# Insert a new blank column before column N ("Late"), shifting the
# existing N/O/P columns (Late, Outstanding heading, Outstanding value)
# one place to the right.
$repayment.Columns("N").Insert()

# The newly inserted column picks up the same display width as column M
# (11 characters, same as Excel copying the width of the column to its
# left on insert) but without the "best fit" auto-sizing flag.
$repayment.Columns("N").ColumnWidth = 10.17

# Update the saved selection / active-tab bookmarks: "Edit Repayment
# Schedule" is no longer the active tab, and its remembered selection
# moves to C7.
$editSchedule.Activate()
$editSchedule.Range("C7").Select() | Out-Null

# "Repayment schedule" becomes the active tab, with the selection resting
# on the newly-inserted column at row 18.
$repayment.Activate()
$repayment.Range("N18").Select() | Out-Null
